$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 82: fix the date/time value in column A ---
$ws.Range("A82").Value = 45461.2916666667

# --- Append new row 83 with the new OHLC data point ---
$ws.Range("A83").Value = 45462.4266550926

# Reuse A82's date/time style (same numFmtId/font) on A83 instead of
# creating a brand new style entry.
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B83").Value = 1500
$ws.Range("C83").Value = 2.99000000953674
$ws.Range("D83").Value = 2.99000000953674
$ws.Range("E83").Value = 2.99000000953674
$ws.Range("F83").Value = 2.99000000953674

# G83 (adj_close) is stored as text in this sheet, like the other rows.
# Force text entry, then clear the resulting formatting so no stray
# number-format style is left behind on the cell (matches sibling cells).
$ws.Range("G83").NumberFormat = "@"
$ws.Range("G83").Value = "2.99000000953674"
$ws.Range("G83").ClearFormats()

# H83 (ticker) is plain text and already detected as such.
$ws.Range("H83").Value = "ESPE.MI"
